$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Fruta / hortaliza, semanal" - rows 2, 4 and 5 (Corazón de apio records for
# Agrícola del Norte S.A. de Arica) get their weekly observations realigned:
# the data that used to sit in row 5 moves to row 2, the data that used to
# sit in row 2 moves to row 4, and the data that used to sit in row 4 moves
# to row 5 (dates D and the associated Calidad/Volumen/Precio columns).

# Row 2 <- previous row 5 values
$ws.Range("D2").Value = 44377
$ws.Range("I2").Value = "Segunda"
$ws.Range("J2").Value = 550
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 2800
$ws.Range("M2").Value = 2364
$ws.Range("P2").Value = 394

# Row 4 <- previous row 2 values
$ws.Range("D4").Value = 44267
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 1500
$ws.Range("L4").Value = 1800
$ws.Range("M4").Value = 1650
$ws.Range("P4").Value = 275

# Row 5 <- previous row 4 values
$ws.Range("D5").Value = 44370
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 1000
$ws.Range("L5").Value = 1200
$ws.Range("M5").Value = 1080
$ws.Range("P5").Value = 180
